$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "M1"
$ws.Range("M2").Value = "M3"
$ws.Range("P2").Value = "M1"
$ws.Range("Q2").Value = "DO"
$ws.Range("R2").Value = "M3"
$ws.Range("T2").Value = "M1"
$ws.Range("W2").Value = "M1"
$ws.Range("AC2").Value = "DO"
$ws.Range("C3").Value = "M3"
$ws.Range("D3").Value = "M1"
$ws.Range("F3").Value = "M1"
$ws.Range("G3").Value = "A2"
$ws.Range("H3").Value = "M1"
$ws.Range("I3").Value = "M1"
$ws.Range("K3").Value = "M2"
$ws.Range("L3").Value = "M3"
$ws.Range("M3").Value = "DO"
$ws.Range("N3").Value = "A2"
$ws.Range("O3").Value = "M1"
$ws.Range("Q3").Value = "M1"
$ws.Range("S3").Value = "M1"
$ws.Range("T3").Value = "M3"
$ws.Range("U3").Value = "DO"
$ws.Range("V3").Value = "M2"
$ws.Range("W3").Value = "M3"
$ws.Range("X3").Value = "M1"
$ws.Range("Y3").Value = "M1"
$ws.Range("AA3").Value = "DO"
$ws.Range("AB3").Value = "A2"
$ws.Range("AC3").Value = "M1"
$ws.Range("B4").Value = "M1"
$ws.Range("C4").Value = "DO"
$ws.Range("F4").Value = "M1"
$ws.Range("H4").Value = "M3"
$ws.Range("K4").Value = "M1"
$ws.Range("L4").Value = "DO"
$ws.Range("P4").Value = "M3"
$ws.Range("S4").Value = "M1"
$ws.Range("T4").Value = "DO"
$ws.Range("X4").Value = "M1"
$ws.Range("Y4").Value = "DO"
$ws.Range("AA4").Value = "M1"
$ws.Range("B5").Value = "DO"
$ws.Range("C5").Value = "M1"
$ws.Range("G5").Value = "M2"
$ws.Range("J5").Value = "M2"
$ws.Range("O5").Value = "A2"
$ws.Range("P5").Value = "M1"
$ws.Range("Q5").Value = "M2"
$ws.Range("R5").Value = "M1"
$ws.Range("S5").Value = "M3"
$ws.Range("T5").Value = "M2"
$ws.Range("U5").Value = "DO"
$ws.Range("V5").Value = "A1"
$ws.Range("X5").Value = "DO"
$ws.Range("Y5").Value = "M1"
$ws.Range("AA5").Value = "M2"
$ws.Range("B6").Value = "DO"
$ws.Range("C6").Value = "M3"
$ws.Range("D6").Value = "A2"
$ws.Range("F6").Value = "A2"
$ws.Range("G6").Value = "M2"
$ws.Range("H6").Value = "A2"
$ws.Range("J6").Value = "A1"
$ws.Range("K6").Value = "A1"
$ws.Range("M6").Value = "A1"
$ws.Range("P6").Value = "A1"
$ws.Range("Q6").Value = "A2"
$ws.Range("S6").Value = "DO"
$ws.Range("T6").Value = "A2"
$ws.Range("X6").Value = "A2"
$ws.Range("Y6").Value = "A1"
$ws.Range("Z6").Value = "A2"
$ws.Range("AA6").Value = "A1"
$ws.Range("AC6").Value = "A2"
$ws.Range("B7").Value = "M3"
$ws.Range("C7").Value = "A1"
$ws.Range("D7").Value = "DO"
$ws.Range("I7").Value = "M3"
$ws.Range("J7").Value = "A1"
$ws.Range("K7").Value = "DO"
$ws.Range("L7").Value = "A1"
$ws.Range("N7").Value = "A1"
$ws.Range("O7").Value = "A1"
$ws.Range("R7").Value = "A1"
$ws.Range("U7").Value = "DO"
$ws.Range("W7").Value = "M3"
$ws.Range("X7").Value = "A1"
$ws.Range("C8").Value = "A1"
$ws.Range("D8").Value = "A2"
$ws.Range("E8").Value = "A2"
$ws.Range("F8").Value = "A1"
$ws.Range("G8").Value = "M2"
$ws.Range("I8").Value = "M3"
$ws.Range("K8").Value = "A2"
$ws.Range("L8").Value = "A2"
$ws.Range("M8").Value = "A2"
$ws.Range("N8").Value = "DO"
$ws.Range("O8").Value = "M1"
$ws.Range("Q8").Value = "A2"
$ws.Range("R8").Value = "A1"
$ws.Range("S8").Value = "A2"
$ws.Range("T8").Value = "A1"
$ws.Range("U8").Value = "DO"
$ws.Range("V8").Value = "M1"
$ws.Range("W8").Value = "M3"
$ws.Range("Y8").Value = "DO"
$ws.Range("Z8").Value = "A1"
$ws.Range("AA8").Value = "A2"
$ws.Range("AB8").Value = "M1"
$ws.Range("AC8").Value = "M1"
$ws.Range("B9").Value = "M1"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = "M3"
$ws.Range("F9").Value = "M2"
$ws.Range("G9").Value = "A1"
$ws.Range("H9").Value = "A2"
$ws.Range("I9").Value = "DO"
$ws.Range("K9").Value = "M2"
$ws.Range("L9").Value = "M3"
$ws.Range("M9").Value = "M2"
$ws.Range("N9").Value = "A2"
$ws.Range("O9").Value = "A2"
$ws.Range("Q9").Value = "M1"
$ws.Range("R9").Value = "M2"
$ws.Range("U9").Value = "A1"
$ws.Range("V9").Value = "A2"
$ws.Range("W9").Value = "M3"
$ws.Range("Y9").Value = "DO"
$ws.Range("Z9").Value = "M1"
$ws.Range("AA9").Value = "M2"
$ws.Range("AB9").Value = "A1"
$ws.Range("B10").Value = "M2"
$ws.Range("E10").Value = "DO"
$ws.Range("G10").Value = "A1"
$ws.Range("H10").Value = "A1"
$ws.Range("I10").Value = "M3"
$ws.Range("J10").Value = "M1"
$ws.Range("K10").Value = "M2"
$ws.Range("L10").Value = "DO"
$ws.Range("M10").Value = "M2"
$ws.Range("N10").Value = "M1"
$ws.Range("O10").Value = "A2"
$ws.Range("Q10").Value = "M1"
$ws.Range("S10").Value = "M1"
$ws.Range("T10").Value = "DO"
$ws.Range("U10").Value = "M2"
$ws.Range("V10").Value = "M1"
$ws.Range("W10").Value = "M1"
$ws.Range("X10").Value = "M2"
$ws.Range("Y10").Value = "M2"
$ws.Range("Z10").Value = "M2"
$ws.Range("AA10").Value = "M2"
$ws.Range("AB10").Value = "M3"
$ws.Range("AC10").Value = "M3"
$ws.Range("R10").Value = "M2"
